$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure changed cells keep their existing Text storage (values are
# percentages/decimal strings stored as text in the source sheet),
# matching the original inlineStr cell type rather than being
# auto-converted to numbers by Excel.
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "332.83"
$ws.Range("E2").Value = "1.65%"
$ws.Range("D3").Value = "45.87"
$ws.Range("E3").Value = "4.37%"
$ws.Range("D4").Value = "5.656"
$ws.Range("E4").Value = "2.76%"
$ws.Range("D5").Value = "0.08385"
$ws.Range("E5").Value = "4.68%"
$ws.Range("D6").Value = "2.037"
$ws.Range("E6").Value = "1.63%"
$ws.Range("D7").Value = "0.9944"
$ws.Range("E7").Value = "4.76%"
$ws.Range("E8").Value = "-0.22%"
$ws.Range("D9").Value = "0.1151"
$ws.Range("E9").Value = "2.38%"
$ws.Range("D10").Value = "0.1928"
$ws.Range("E10").Value = "2.95%"
$ws.Range("E11").Value = "-2.78%"
$ws.Range("D12").Value = "0.09971"
$ws.Range("E12").Value = "1.25%"
$ws.Range("D13").Value = "0.04683"
$ws.Range("E13").Value = "2.17%"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").Value = "-0.71%"
$ws.Range("D15").Value = "0.001282"
$ws.Range("E15").Value = "1.32%"
$ws.Range("D16").Value = "0.006097"
$ws.Range("E16").Value = "3.02%"
$ws.Range("D17").Value = "3.376"
$ws.Range("E17").Value = "0.66%"
$ws.Range("D18").Value = "4.487"
$ws.Range("E18").Value = "3.92%"
$ws.Range("D19").Value = "0.3365"
$ws.Range("E19").Value = "-3.20%"
$ws.Range("D20").Value = "0.1403"
$ws.Range("E20").Value = "-0.21%"
$ws.Range("D21").Value = "0.2656"
$ws.Range("E21").Value = "4.46%"
$ws.Range("E22").Value = "3.70%"
$ws.Range("D23").Value = "0.001311"
$ws.Range("E23").Value = "4.27%"
$ws.Range("D24").Value = "0.004637"
$ws.Range("E24").Value = "7.08%"
$ws.Range("E25").Value = "10.86%"
$ws.Range("D26").Value = "0.0003753"
$ws.Range("E26").Value = "0.36%"
$ws.Range("D38").Value = "0.02779"
$ws.Range("E38").Value = "8.30%"
$ws.Range("D39").Value = "0.05747"
$ws.Range("E39").Value = "1.15%"
$ws.Range("D40").Value = "0.007757"
$ws.Range("E40").Value = "2.75%"
$ws.Range("D41").Value = "0.1435"
$ws.Range("E41").Value = "2.70%"
$ws.Range("D42").Value = "0.007288"
$ws.Range("E42").Value = "-4.20%"
$ws.Range("D43").Value = "0.002114"
$ws.Range("E43").Value = "5.04%"
$ws.Range("D44").Value = "0.009047"
$ws.Range("E44").Value = "1.96%"
$ws.Range("D45").Value = "0.3412"
$ws.Range("D46").Value = "0.00007367"
$ws.Range("E46").Value = "3.81%"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").Value = "0.47%"
$ws.Range("D48").Value = "0.0005816"
$ws.Range("E48").Value = "0.09%"
$ws.Range("D49").Value = "0.003512"
$ws.Range("E49").Value = "0.88%"
$ws.Range("D50").Value = "0.003510"
$ws.Range("E50").Value = "-0.54%"
$ws.Range("D51").Value = "0.00002107"
$ws.Range("E51").Value = "0.47%"
